$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) hold numeric-looking text that Excel would otherwise
# auto-convert to a Number when assigned via .Value. Force these specific
# cells to Text format first so the values persist exactly as strings,
# matching the original inlineStr cell content.

$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.317.96"
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").Value = "1.930.41"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "248.76"
$ws.Range("E5").Value = "  -4.01%  "
$ws.Range("D6").Value = "0.7178"
$ws.Range("E6").Value = "  -8.96%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.3293"
$ws.Range("E8").Value = "  -8.13%  "
$ws.Range("D9").Value = "27.91"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "0.06915"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -6.06%  "
$ws.Range("D12").Value = "0.08058"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "1.931.35"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "5.397"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").Value = "94.66"
$ws.Range("E15").Value = "  -6.32%  "
$ws.Range("D16").Value = "14.48"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.304.50"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008326"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("D19").Value = "252.84"
$ws.Range("E19").Value = "  -8.80%  "
$ws.Range("D20").Value = "5.804"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "2.182.66"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "6.860"
$ws.Range("E24").Value = "  -4.97%  "
$ws.Range("D25").Value = "9.705"
$ws.Range("E25").Value = "  -3.76%  "
$ws.Range("D26").Value = "159.39"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").Value = "2.394"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "19.10"
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("D29").Value = "0.1334"
$ws.Range("E29").Value = "  -11.85%  "
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").Value = "4.397"
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").Value = "4.185"
$ws.Range("E33").Value = "  -5.27%  "
$ws.Range("D34").Value = "0.05107"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").Value = "1.219"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "0.7391"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("D37").Value = "2.744"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").Value = "2.833"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "6.591"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "78.74"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  -6.53%  "
$ws.Range("D43").Value = "1.988"
$ws.Range("E43").Value = "  -8.49%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "0.8362"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("D46").Value = "101.83"
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").Value = "9.769"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "7.291"
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("D49").Value = "36.49"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").Value = "0.05954"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "0.4072"
